# MR Hassani Degree Added with failed
# Row 12 (MohammadReza Hassani) gets its quiz scores filled in, matching the
# same "flagged" (orange) formatting already used for other students with an
# asterisk / special-case score (e.g. row 6 / row 8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing "flag" formatting (orange fill + fonts/borders) from
# row 6 columns B:K onto row 12 columns B:K so the new scores get the same
# look-and-feel as other already-graded rows that carry a "*" footnote.
$ws.Range("B6:K6").Copy()
$ws.Range("B12:K12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Highlight the student's name/email cell (A12) with the same orange fill
# used for the rest of the row.
$ws.Range("A12").Interior.Color = 49407

# Fill in MohammadReza Hassani's quiz grades.
$ws.Range("F12").Value = 100
$ws.Range("G12").Value = 95
$ws.Range("H12").Value = "*"
$ws.Range("I12").Value = 93
$ws.Range("J12").Value = 85
$ws.Range("K12").Value = 100

# Move the active selection to the last entered cell and reset the
# scrolled/frozen top-left cell back to the top of the sheet.
$ws.Application.ActiveWindow.ScrollRow = 1
$null = $ws.Range("J12").Select()
